$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-numeric-looking text cells (Coin name, Link, Volume%) - safe to assign directly.
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("E24").Value = "  -4.24%  "
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("E36").Value = "  +8.81%  "
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("E44").Value = "  +3.25%  "
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  -3.04%  "

# Price column (D) cells: many values parse as valid numbers (e.g. "1.00", "65.16").
# The source data stores these as plain text, so force text entry via a temporary
# "@" (Text) number format, then ClearFormats() to drop the temporary style again
# (re-using cellXfs index 0, so no stray style/format diff remains).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.886.98"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.25"
$ws.Range("D3").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.71"
$ws.Range("D5").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.21"
$ws.Range("D8").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.863.68"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.23"
$ws.Range("D13").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.566"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.16"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.886.97"
$ws.Range("D17").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0721"
$ws.Range("D19").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.33"
$ws.Range("D23").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.62"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("D26").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.60"
$ws.Range("D28").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.06"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.395.60"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0170"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.559"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.869"
$ws.Range("D40").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.82"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.52"
$ws.Range("D44").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.773.44"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.64"
$ws.Range("D48").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0506"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.48"
$ws.Range("D51").ClearFormats()
